$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-9:
# value changes from 45170 (2023-09-01) to 45174 (2023-09-05)
$newDate = [DateTime]::FromOADate(45174)

foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
